# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force plain text storage so numeric-looking strings (e.g. prices like
    # '579.24') are not silently reinterpreted as numbers by Excel.
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range('D2') '60.375.61'
Set-TextValue $ws.Range('E2') '  -4.39%  '
Set-TextValue $ws.Range('D3') '2.997.56'
Set-TextValue $ws.Range('E3') '  -5.76%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '579.24'
Set-TextValue $ws.Range('E5') '  -2.14%  '
Set-TextValue $ws.Range('D6') '126.58'
Set-TextValue $ws.Range('E6') '  -6.41%  '
Set-TextValue $ws.Range('E7') '  +0.10%  '
Set-TextValue $ws.Range('D8') '2.993.37'
Set-TextValue $ws.Range('E8') '  -5.78%  '
Set-TextValue $ws.Range('E9') '  -3.13%  '
Set-TextValue $ws.Range('D10') '0.133'
Set-TextValue $ws.Range('E10') '  -5.92%  '
Set-TextValue $ws.Range('D11') '5.13'
Set-TextValue $ws.Range('E11') '  -2.16%  '
Set-TextValue $ws.Range('E12') '  -3.11%  '
Set-TextValue $ws.Range('D13') '0.0000224'
Set-TextValue $ws.Range('E13') '  -5.63%  '
Set-TextValue $ws.Range('D14') '32.56'
Set-TextValue $ws.Range('E14') '  -5.59%  '
Set-TextValue $ws.Range('E15') '  +0.27%  '
Set-TextValue $ws.Range('D16') '3.485.41'
Set-TextValue $ws.Range('E16') '  -5.91%  '
Set-TextValue $ws.Range('D17') '2.992.33'
Set-TextValue $ws.Range('E17') '  -5.90%  '
Set-TextValue $ws.Range('D18') '60.230.01'
Set-TextValue $ws.Range('E18') '  -4.63%  '
Set-TextValue $ws.Range('D19') '6.29'
Set-TextValue $ws.Range('E19') '  -3.99%  '
Set-TextValue $ws.Range('D20') '431.80'
Set-TextValue $ws.Range('E20') '  -6.31%  '
Set-TextValue $ws.Range('D21') '13.10'
Set-TextValue $ws.Range('E21') '  -6.08%  '
Set-TextValue $ws.Range('D22') '0.662'
Set-TextValue $ws.Range('E22') '  -4.91%  '
Set-TextValue $ws.Range('D23') '7.05'
Set-TextValue $ws.Range('E23') '  -7.47%  '
Set-TextValue $ws.Range('D24') '12.80'
Set-TextValue $ws.Range('E24') '  -3.92%  '
Set-TextValue $ws.Range('D25') '79.28'
Set-TextValue $ws.Range('E25') '  -3.85%  '
Set-TextValue $ws.Range('E26') '  +0.09%  '
Set-TextValue $ws.Range('D27') '0.998'
Set-TextValue $ws.Range('E27') '  -0.18%  '
Set-TextValue $ws.Range('E28') '  -4.13%  '
Set-TextValue $ws.Range('D29') '7.29'
Set-TextValue $ws.Range('E29') '  -5.22%  '
Set-TextValue $ws.Range('D30') '1.89'
Set-TextValue $ws.Range('E30') '  -7.03%  '
Set-TextValue $ws.Range('D31') '6.15'
Set-TextValue $ws.Range('E31') '  -9.01%  '
Set-TextValue $ws.Range('D32') '25.39'
Set-TextValue $ws.Range('E32') '  -6.95%  '
Set-TextValue $ws.Range('D33') '0.0940'
Set-TextValue $ws.Range('E33') '  -7.80%  '
Set-TextValue $ws.Range('D34') '2.15'
Set-TextValue $ws.Range('E34') '  -9.37%  '
Set-TextValue $ws.Range('E35') '  -7.23%  '
Set-TextValue $ws.Range('D36') '5.60'
Set-TextValue $ws.Range('D37') '49.73'
Set-TextValue $ws.Range('E37') '  -2.88%  '
Set-TextValue $ws.Range('D38') '0.0₃0662'
Set-TextValue $ws.Range('E38') '  -7.16%  '
Set-TextValue $ws.Range('D39') '8.23'
Set-TextValue $ws.Range('E39') '  +1.72%  '
Set-TextValue $ws.Range('D40') '0.0360'
Set-TextValue $ws.Range('E41') '  -1.13%  '
Set-TextValue $ws.Range('D42') '382.80'
Set-TextValue $ws.Range('E42') '  -4.99%  '
Set-TextValue $ws.Range('D43') '2.47'
Set-TextValue $ws.Range('E43') '  -6.98%  '
Set-TextValue $ws.Range('D44') '2.649.03'
Set-TextValue $ws.Range('E44') '  -5.63%  '
Set-TextValue $ws.Range('E46') '  -6.34%  '
Set-TextValue $ws.Range('D47') '2.00'
Set-TextValue $ws.Range('E47') '  -5.70%  '
Set-TextValue $ws.Range('D48') '118.51'
Set-TextValue $ws.Range('E48') '  -6.23%  '
Set-TextValue $ws.Range('E49') '  -3.95%  '
Set-TextValue $ws.Range('D50') '23.69'
Set-TextValue $ws.Range('E50') '  -5.98%  '
Set-TextValue $ws.Range('B51') 'Cronos'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.133'
Set-TextValue $ws.Range('E51') '  +4.04%  '
